# "laboratorio 7 - Entrega Final"
#
# The lab's placeholder/demo measurements in the "Datos Lab7" sheet are
# replaced with the real measured results (execution time vs. memory
# consumption for PROBING and CHAINING collision-resolution strategies),
# and the active selection is moved to where the student left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab7")

# --- Tabla "Factor de Carga (PROBING)" (filas 3 a 5) ---------------------
# Columna B: Consumo de Datos [kB]   Columna C: Tiempo de Ejecución [ms]
$ws.Range("B3").Value = 275591
$ws.Range("C3").Value = 4628444

$ws.Range("B4").Value = 238936
$ws.Range("C4").Value = 5833070

$ws.Range("B5").Value = 275591
$ws.Range("C5").Value = 7326676

# --- Tabla "Factor de Carga (CHAINING)" (filas 10 a 12) -------------------
$ws.Range("B10").Value = 247385
$ws.Range("C10").Value = 6339687

$ws.Range("B11").Value = 241840
$ws.Range("C11").Value = 6454886

$ws.Range("B12").Value = 240152
$ws.Range("C12").Value = 6111928

# Deja la selección en C12, tal como quedó la hoja al guardar la entrega.
$ws.Range("C12").Select()
